$d = $word.ActiveDocument

# --- Paragraph 1: "This is a Microsoft word document." ---------------------
$p1 = $d.Paragraphs(1)
$r = $p1.Range
[void]$r.MoveEnd(1, -1)    # exclude the paragraph mark from the range
[void]$r.Collapse(0)       # collapse to the end of the visible text

# Two trailing spaces appended to the existing (uncolored) run
[void]$r.InsertAfter("  ")
[void]$r.Collapse(0)

# Insert the new, completely empty paragraph right after paragraph 1, before
# any colored text is added, so the blank paragraph doesn't inherit red
# character formatting.
[void]$r.InsertParagraphAfter()

# Re-acquire paragraph 1's range (it still ends right before the new blank
# paragraph that was just inserted).
$p1 = $d.Paragraphs(1)
$r = $p1.Range
[void]$r.MoveEnd(1, -1)
[void]$r.Collapse(0)

# Run: red "(This is a change – Version for branch "
$rBranchLabel = $r.Duplicate
[void]$rBranchLabel.Collapse(0)
[void]$rBranchLabel.InsertAfter("(This is a change " + [char]0x2013 + " Version for branch ")
$rBranchLabel.Font.Color = 192  # RGB(0xC0,0x00,0x00) -> w:color C00000

# Run: red branch name "main" (its own run, separate from the label text)
$rBranchName = $rBranchLabel.Duplicate
[void]$rBranchName.Collapse(0)
[void]$rBranchName.InsertAfter("main")
$rBranchName.Font.Color = 192

# Run: red closing paren ")"
$rClose = $rBranchName.Duplicate
[void]$rClose.Collapse(0)
[void]$rClose.InsertAfter(")")
$rClose.Font.Color = 192
